$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPRT")

# Row 4 (Inventory)
$ws.Range("B4").Value = 30000000.0
$ws.Range("C4").Value = 28000000.0
$ws.Range("D4").Value = 20000000.0
$ws.Range("E4").Value = 20000000.0
$ws.Range("F4").Value = 19000000.0

# Row 14 (Accounts Payable)
$ws.Range("B14").Value = 324000000.0
$ws.Range("C14").Value = 360000000.0
$ws.Range("D14").Value = 133000000.0
$ws.Range("E14").Value = 251000000.0
$ws.Range("F14").Value = 297000000.0

# Row 23 (Long Term Tax Liability (Deferred))
$ws.Range("B23").Value = 81000000.0
$ws.Range("C23").Value = 78000000.0
$ws.Range("D23").Value = 71000000.0
$ws.Range("E23").Value = 57000000.0
$ws.Range("F23").Value = 55000000.0
